$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "44.065.00"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.65%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.360.32"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.60%  "

$ws.Range("E4").Value = "  -0.14%  "

$ws.Range("E5").Value = "  +3.88%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "238.46"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.21%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "72.84"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +11.47%  "

$ws.Range("E8").Value = "  -0.05%  "

$ws.Range("E9").Value = "  +18.01%  "

$ws.Range("E10").Value = "  +6.01%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "29.31"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +9.95%  "

$ws.Range("E12").Value = "  +2.83%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.713.28"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.75%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "16.83"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +9.77%  "

$ws.Range("E15").Value = "  +7.10%  "

$ws.Range("E16").Value = "  +8.21%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.367.73"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.92%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "44.049.75"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.78%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0000102"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.94%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "78.03"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +6.06%  "

$ws.Range("E21").Value = "  +4.62%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "255.35"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.89%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.76"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.83%  "

$ws.Range("E25").Value = "  +3.32%  "

$ws.Range("E26").Value = "  +6.52%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.24"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.20%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "22.45"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.16%  "

$ws.Range("E29").Value = "  +5.49%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "172.79"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.53%  "

$ws.Range("E31").Value = "  +2.51%  "

$ws.Range("E32").Value = "  +5.21%  "

$ws.Range("E33").Value = "  +4.39%  "

$ws.Range("E34").Value = "  +6.69%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.28"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.84%  "

$ws.Range("E36").Value = "  +10.34%  "

$ws.Range("E37").Value = "  -1.32%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.45"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.10%  "

$ws.Range("E39").Value = "  +7.41%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "19.59"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +9.58%  "

$ws.Range("E41").Value = "  +0.01%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.85"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.32%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.25"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.44%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0983"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.02%  "

$ws.Range("E45").Value = "  +1.44%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "98.64"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.20%  "

$ws.Range("E47").Value = "  +1.54%  "

$ws.Range("E48").Value = "  +12.25%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.35"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.99%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.440.62"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.39%  "

$ws.Range("E51").Value = "  +1.47%  "
